$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in row 2 (E2, F2) and match D2's number format/style
$ws.Range("E2").Value = 42844
$ws.Range("F2").Value = 42845
$ws.Range("E2").NumberFormat = "DD/MM/YY"
$ws.Range("F2").NumberFormat = "DD/MM/YY"

# Widen column E to fit the new structure (closest attainable width to 10.25 chars
# given the host's pixel-grid rounding on the ColumnWidth COM property)
$ws.Columns.Item(5).ColumnWidth = 9.42

# Update the view selection for the bottom-right frozen pane to F14
$ws.Range("F14").Select()
